$wb = $excel.ActiveWorkbook

$newTime = "01:58:07"

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"
$ws1.Range("A3").Value = "Total filas: 3"

# Existing row 6 (14_ABASTO) - scrape time + arrival + minutes updated
$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "02:01"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 3
$ws1.Range("E6").Value = "LP1912"

# Existing row 7 (215_ALUAR) - scrape time + arrival + minutes updated
$ws1.Range("A7").Value = $newTime
$ws1.Range("B7").Value = "03:02"
$ws1.Range("C7").Value = "215_ALUAR"
$ws1.Range("D7").Value = 64
$ws1.Range("E7").Value = "LP1912"

# New row 8 (14_ABASTO)
$ws1.Range("A8").Value = $newTime
$ws1.Range("B8").Value = "03:51"
$ws1.Range("C8").Value = "14_ABASTO"
$ws1.Range("D8").Value = 113
$ws1.Range("E8").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "03:02"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 64
$ws2.Range("E6").Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
